{"js": "// Merge the \"Table S4. \" caption runs (\"Table S\" + \"4\" + \". \") into a\n// single run reading \"Table S. \" (the table number is dropped), matching\n// the target OOXML diff. We search the body for the exact text (search()\n// transparently spans run boundaries) and replace the hit in place so the\n// surrounding bold/font formatting of the first run is preserved.\nconst body = context.document.body;\nconst results = body.search(\"Table S4. \", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Table S4. \" in the document body.');\n}\n\nresults.items[0].insertText(\"Table S. \", \"Replace\");\nawait context.sync();\n", "ps1": "# Merge the \"Table S4. \" caption runs (\"Table S\" + \"4\" + \". \") into a\n# single run reading \"Table S. \" (the table number is dropped), matching\n# the target OOXML diff. Find/Replace spans run boundaries and Word will\n# fold the whole match into the formatting of the first run, collapsing\n# the now-redundant extra runs - exactly mirroring the diff.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Table S4. \"\n$find.Replacement.Text = \"Table S. \"\n\n# FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n# MatchAllWordForms, Forward, Wrap (1=wdFindContinue), Format, ReplaceWith,\n# Replace (2=wdReplaceAll)\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
